{"js": "// The document ends with an empty paragraph right after the\n// \"GitHub Link:\" heading. Turn that empty paragraph into a run holding\n// the repository URL, styled blue (#0070C0) at 14pt (sz/szCs = 28\n// half-points), matching the target OOXML exactly.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// The target is the very last paragraph of the body (immediately before\n// the sectPr) -- in this document that is also the empty one.\nlet targetParagraph = paragraphs.items[paragraphs.items.length - 1];\nif (targetParagraph.text !== \"\") {\n  // Defensive fallback: locate the trailing empty paragraph explicitly\n  // in case the body ever gains more paragraphs after it.\n  for (let i = paragraphs.items.length - 1; i >= 0; i--) {\n    if (paragraphs.items[i].text === \"\") {\n      targetParagraph = paragraphs.items[i];\n      break;\n    }\n  }\n}\n\nconst url = \"https://github.com/Vasanth30e/Practice_Project_Phase3/tree/master/Practice_Project%203\";\n\n// Word.Font (color/size/bold/...) has no setter for the complex-script\n// size (w:szCs), so drive this edit through insertOoxml: it lets us\n// write the paragraph-mark run properties (w:pPr/w:rPr) and the run's\n// own properties (w:r/w:rPr) exactly as authored, including w:szCs.\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n      '<pkg:xmlData>' +\n        '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n          '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n        '</Relationships>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body>' +\n            '<w:p>' +\n              '<w:pPr>' +\n                '<w:rPr>' +\n                  '<w:color w:val=\"0070C0\"/>' +\n                  '<w:sz w:val=\"28\"/>' +\n                  '<w:szCs w:val=\"28\"/>' +\n                '</w:rPr>' +\n              '</w:pPr>' +\n              '<w:r>' +\n                '<w:rPr>' +\n                  '<w:color w:val=\"0070C0\"/>' +\n                  '<w:sz w:val=\"28\"/>' +\n                  '<w:szCs w:val=\"28\"/>' +\n                '</w:rPr>' +\n                '<w:t>' + url + '</w:t>' +\n              '</w:r>' +\n            '</w:p>' +\n          '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\ntargetParagraph.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The document ends with an empty paragraph right after the\n# \"GitHub Link:\" heading. Fill that empty paragraph with the repository\n# URL, styled blue (#0070C0) at 14pt, matching the target OOXML exactly.\n$d = $word.ActiveDocument\n\n# Target is the last paragraph in the document (immediately before the\n# sectPr) -- defensively confirm it is empty, and otherwise fall back to\n# scanning backwards for the trailing empty paragraph.\n$count = $d.Paragraphs.Count\n$p = $d.Paragraphs.Item($count)\n\nif ($p.Range.Text -ne \"\") {\n    for ($i = $count; $i -ge 1; $i--) {\n        $candidate = $d.Paragraphs.Item($i)\n        if ($candidate.Range.Text -eq \"\") {\n            $p = $candidate\n            break\n        }\n    }\n}\n\n$r = $p.Range\n$r.Text = \"https://github.com/Vasanth30e/Practice_Project_Phase3/tree/master/Practice_Project%203\"\n\n# wdColor values are packed as 0x00BBGGRR, so RGB 0070C0 -> BGR C07000.\n$r.Font.Color = 12611584\n$r.Font.Size = 14\n# Also set the complex-script size so w:szCs is written alongside w:sz.\n$r.Font.SizeBi = 14\n"}
